$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "F2"
$ws.Range("C2").Value = "F2rl1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.4824693333333334
$ws.Range("H2").Value = 1.447408
$ws.Range("I2").Value = 0.2099856823459656
$ws.Range("J2").Value = 0.2099856823459655
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.295943
$ws.Range("N2").Value = 0.887829
$ws.Range("O2").Value = 0.02818120871217195
$ws.Range("P2").Value = 0.02818120871217194
$ws.Range("Q2").Value = 0.1427834219146667
$ws.Range("R2").Value = 1.285050797232
$ws.Range("S2").Value = 0.005917650340759495
$ws.Range("T2").Value = 0.005917650340759493

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "F2"
$ws.Range("C3").Value = "F2rl1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.4824693333333334
$ws.Range("H3").Value = 1.447408
$ws.Range("I3").Value = 0.2099856823459656
$ws.Range("J3").Value = 0.2099856823459655
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 10.20548733333333
$ws.Range("N3").Value = 30.616462
$ws.Range("O3").Value = 0.9718187912878281
$ws.Range("P3").Value = 0.9718187912878281
$ws.Range("Q3").Value = 4.923834670055111
$ws.Range("R3").Value = 44.314512030496
$ws.Range("S3").Value = 0.2040680320052061
$ws.Range("T3").Value = 0.204068032005206

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "F2"
$ws.Range("C4").Value = "F2rl1"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 1.110828
$ws.Range("H4").Value = 3.332484
$ws.Range("I4").Value = 0.4834669468781523
$ws.Range("J4").Value = 0.4834669468781522
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.295943
$ws.Range("N4").Value = 0.887829
$ws.Range("O4").Value = 0.02818120871217195
$ws.Range("P4").Value = 0.02818120871217194
$ws.Range("Q4").Value = 0.328741770804
$ws.Range("R4").Value = 2.958675937236
$ws.Range("S4").Value = 0.01362468293540976
$ws.Range("T4").Value = 0.01362468293540975

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "F2"
$ws.Range("C5").Value = "F2rl1"
$ws.Range("D5").Value = "FAPs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 1.110828
$ws.Range("H5").Value = 3.332484
$ws.Range("I5").Value = 0.4834669468781523
$ws.Range("J5").Value = 0.4834669468781522
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 10.20548733333333
$ws.Range("N5").Value = 30.616462
$ws.Range("O5").Value = 0.9718187912878281
$ws.Range("P5").Value = 0.9718187912878281
$ws.Range("Q5").Value = 11.336541083512
$ws.Range("R5").Value = 102.028869751608
$ws.Range("S5").Value = 0.4698422639427425
$ws.Range("T5").Value = 0.4698422639427425

# Row 6
$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "F2"
$ws.Range("C6").Value = "F2rl1"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.7043323333333333
$ws.Range("H6").Value = 2.112997
$ws.Range("I6").Value = 0.3065473707758822
$ws.Range("J6").Value = 0.3065473707758822
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.295943
$ws.Range("N6").Value = 0.887829
$ws.Range("O6").Value = 0.02818120871217195
$ws.Range("P6").Value = 0.02818120871217194
$ws.Range("Q6").Value = 0.2084422237236667
$ws.Range("R6").Value = 1.875980013513
$ws.Range("S6").Value = 0.008638875436002696
$ws.Range("T6").Value = 0.008638875436002695

# Row 7
$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "F2"
$ws.Range("C7").Value = "F2rl1"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.7043323333333333
$ws.Range("H7").Value = 2.112997
$ws.Range("I7").Value = 0.3065473707758822
$ws.Range("J7").Value = 0.3065473707758822
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 10.20548733333333
$ws.Range("N7").Value = 30.616462
$ws.Range("O7").Value = 0.9718187912878281
$ws.Range("P7").Value = 0.9718187912878281
$ws.Range("Q7").Value = 7.188054706290444
$ws.Range("R7").Value = 64.692492356614
$ws.Range("S7").Value = 0.2979084953398796
$ws.Range("T7").Value = 0.2979084953398796
